# Insert a new data row at row 852 (shifting existing rows 852:905 down to 853:906)
# and populate it with the new Kiwi price record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 852; Excel shifts rows 852:905 down to 853:906
# and copies formatting (incl. the date style on column D) from the row above.
$ws.Rows.Item(852).Insert()

# Populate the newly inserted row 852 with the new record's values.
$ws.Range("A852").Value2 = 10
$ws.Range("B852").Value2 = "Vega Modelo de Temuco"
$ws.Range("C852").Value2 = "La Araucanía"
$ws.Range("D852").Value2 = 45265
$ws.Range("E852").Value2 = 9
$ws.Range("F852").Value2 = "Fruta"
$ws.Range("G852").Value2 = 100101
$ws.Range("H852").Value2 = "Berries"
$ws.Range("I852").Value2 = 100101007
$ws.Range("J852").Value2 = "Kiwi"
$ws.Range("K852").Value2 = "Hayward"
$ws.Range("L852").Value2 = "Primera"
$ws.Range("M852").Value2 = 100
$ws.Range("N852").Value2 = 26000
$ws.Range("O852").Value2 = 26000
$ws.Range("P852").Value2 = 26000
$ws.Range("Q852").Value2 = "$/bandeja 18 kilos"
$ws.Range("R852").Value2 = "Región de O'Higgins"
$ws.Range("S852").Value2 = 1444
$ws.Range("T852").Value2 = 18
